$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.730.59'
$ws.Range("E2").Value = '  +3.89%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.862.62'
$ws.Range("E3").Value = '  +2.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9985'
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '274.28'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9988'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5275'
$ws.Range("E7").Value = '  +3.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3386'
$ws.Range("E8").Value = '  -4.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06805'
$ws.Range("E9").Value = '  +1.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.89'
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7948'
$ws.Range("E11").Value = '  -3.90%  '
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.821.08'
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '89.86'
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.128'
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9974'
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.43'
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008002'
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9988'
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.771.83'
$ws.Range("E20").Value = '  +3.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.102.75'
$ws.Range("E21").Value = '  +2.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.716'
$ws.Range("E22").Value = '  -0.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.977'
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.107'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.364'
$ws.Range("E25").Value = '  +5.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.71'
$ws.Range("E26").Value = '  +2.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.658'
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.20'
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.32'
$ws.Range("E29").Value = '  +2.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.334'
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.313'
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08873'
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04917'
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.161'
$ws.Range("E34").Value = '  +2.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7288'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.878'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.227'
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.335'
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01846'
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5104'
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9397'
$ws.Range("E41").Value = '  -3.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '116.20'
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.132'
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.008'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9981'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4412'
$ws.Range("E46").Value = '  -2.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1328'
$ws.Range("E47").Value = '  -3.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.305'
$ws.Range("E48").Value = '  +1.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.10'
$ws.Range("E49").Value = '  -1.04%  '
$ws.Range("E50").Value = '  +1.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.475'
$ws.Range("E51").Value = '  -1.67%  '
